# Applies the three text edits described by the commit diff:
#   1. Merge the many single-word runs that spell out the "Mega Plaza
#      Oriente" discount problem into one contiguous run of text.
#   2. Merge the runs that spell out the "Si trabaja mas de 40 horas..."
#      sentence into one contiguous run of text.
#   3. Rework the "Implementa una Clase..." paragraph: insert a new
#      clause after "un metodo" and change "uso" to "usa".
#
# Word's Find/Replace naturally collapses the matched span (and any runs
# it fully covers) into a single run carrying the formatting of the
# first run in the match, which mirrors the merges seen in the diff.

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "En Mega Plaza Oriente hay ... (¡¡¡Imprime el "
# -----------------------------------------------------------------
$old1 = "En Mega Plaza Oriente hay un 20% de descuento a los clientes cuya" `
    + " compra supere los `$300 pesos ¿Cuál será la cantidad que pagará" `
    + " una persona por su compra? (¡¡¡Imprime el "
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, `
    $true, 1, $false, $old1, 2) | Out-Null

# -----------------------------------------------------------------
# Change 2: "Si trabaja más de 40 horas se le paga ..."
# -----------------------------------------------------------------
$old2 = "Si trabaja más de 40 horas se le paga `$16 por cada una de las" `
    + " primeras 40 horas y `$20 por cada hora extra."
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, `
    $true, 1, $false, $old2, 2) | Out-Null

# -----------------------------------------------------------------
# Change 3: "Implementa una Clase ... un método, uso de getters ..."
#   -> "... un método que realice las acciones de un objeto de la
#        clase, usa de getters ..."
# -----------------------------------------------------------------
$old3 = "un método, uso de "
$new3 = "un método que realice las acciones de un objeto de la clase, usa de "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, `
    $true, 1, $false, $new3, 2) | Out-Null
